$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 20718.182
$ws.Range("I7").Value = 1000
$ws.Range("J7").Value = 22690
$ws.Range("K7").Value = 1000
$ws.Range("L7").Value = 22690
$ws.Range("M7").Value = -888
$ws.Range("N7").Value = -22914

$ws.Range("H13").Value = 30000
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 30000
$ws.Range("K13").Value = 0
$ws.Range("L13").Value = 30000
$ws.Range("M13").Value = $null
$ws.Range("N13").Value = -30338

$ws.Range("H14").Value = 20718.182
$ws.Range("I14").Value = 1000
$ws.Range("J14").Value = 22690
$ws.Range("K14").Value = 1000
$ws.Range("L14").Value = 22690
$ws.Range("M14").Value = -809
$ws.Range("N14").Value = -23072

$ws.Range("I20").Value = 8220
$ws.Range("K20").Value = 8220
$ws.Range("M20").Value = -7990

$ws.Range("I35").Value = 8220
$ws.Range("K35").Value = 8220
$ws.Range("M35").Value = -7841

$ws.Range("H129").Value = 1039.9423
$ws.Range("J129").Value = 1139.2273
$ws.Range("L129").Value = 3417.6819
$ws.Range("N129").Value = -13417.6819

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 122
$ws.Range("I5").Value = 122
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 122
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -10
$ws.Range("N5").Value = $null

$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").Value = $null

$ws.Range("H32").Value = 9011.686
$ws.Range("I32").Value = 6701.873
$ws.Range("J32").Value = 29800
$ws.Range("K32").Value = 6701.873
$ws.Range("L32").Value = 29800
$ws.Range("M32").Value = -6414.873
$ws.Range("N32").Value = -30374

$ws.Range("H122").Value = 1225118.6
$ws.Range("I122").Value = 1428482.1
$ws.Range("J122").Value = 4938
$ws.Range("K122").Value = 4285446.300000001
$ws.Range("L122").Value = 14814
$ws.Range("M122").Value = -4282996.300000001
$ws.Range("N122").Value = -19714

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 122
$ws.Range("I4").Value = 122
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 122
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -7
$ws.Range("N4").Value = $null

$ws.Range("H25").Value = 2654.5
$ws.Range("J25").Value = 7220
$ws.Range("L25").Value = 7220
$ws.Range("N25").Value = -7690

$ws.Range("H26").Value = 15117.75
$ws.Range("I26").Value = 15117.75
$ws.Range("K26").Value = 15117.75
$ws.Range("M26").Value = -14825.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 2762.375
$ws.Range("I19").Value = 625
$ws.Range("K19").Value = 625
$ws.Range("M19").Value = -455

$ws.Range("H24").Value = 2762.375
$ws.Range("I24").Value = 625
$ws.Range("K24").Value = 625
$ws.Range("M24").Value = -455

$ws.Range("H141").Value = 411746.25
$ws.Range("I141").Value = 10296
$ws.Range("J141").Value = 469096.28
$ws.Range("K141").Value = 10296
$ws.Range("L141").Value = 469096.28
$ws.Range("M141").Value = -5116
$ws.Range("N141").Value = -479456.28

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 840.8333
$ws.Range("I17").Value = 50
$ws.Range("K17").Value = 150
$ws.Range("M17").Value = 19

$ws.Range("H19").Value = 19999
$ws.Range("J19").Value = 19999
$ws.Range("L19").Value = 59997
$ws.Range("N19").Value = -60345

$ws.Range("H97").Value = 3846505.8
$ws.Range("I97").Value = 5263506.5
$ws.Range("J97").Value = 360.57144
$ws.Range("K97").Value = 15790519.5
$ws.Range("L97").Value = 1081.71432
$ws.Range("M97").Value = -15790023.5
$ws.Range("N97").Value = -2073.71432

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H4").Value = 40000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 40000
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 40000
$ws.Range("M4").Value = $null
$ws.Range("N4").Value = -40224

$ws.Range("H12").Value = 23144.445
$ws.Range("I12").Value = 24187.5
$ws.Range("J12").Value = 14800
$ws.Range("K12").Value = 24187.5
$ws.Range("L12").Value = 14800
$ws.Range("M12").Value = -24047.5
$ws.Range("N12").Value = -15080

$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = $null

$ws.Range("H123").Value = 18408.97
$ws.Range("J123").Value = 18694.092
$ws.Range("L123").Value = 18694.092
$ws.Range("N123").Value = -23594.092

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2112
$ws.Range("I22").Value = 300
$ws.Range("J22").Value = 2293.2
$ws.Range("K22").Value = 300
$ws.Range("L22").Value = 2293.2
$ws.Range("M22").Value = -5
$ws.Range("N22").Value = -2883.2

$ws.Range("H27").Value = 2112
$ws.Range("I27").Value = 300
$ws.Range("J27").Value = 2293.2
$ws.Range("K27").Value = 300
$ws.Range("L27").Value = 2293.2
$ws.Range("M27").Value = -193
$ws.Range("N27").Value = -2507.2

$ws.Range("H132").Value = 20836914
$ws.Range("I132").Value = 33336984
$ws.Range("J132").Value = 3466.5
$ws.Range("K132").Value = 100010952
$ws.Range("L132").Value = 10399.5
$ws.Range("M132").Value = -100008422
$ws.Range("N132").Value = -15459.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 10000
$ws.Range("J18").Value = 10000
$ws.Range("L18").Value = 10000
$ws.Range("N18").Value = -10346

$ws.Range("H64").Value = 25114
$ws.Range("J64").Value = 25114
$ws.Range("L64").Value = 25114
$ws.Range("N64").Value = -25610

$ws.Range("H67").Value = 25114
$ws.Range("J67").Value = 25114
$ws.Range("L67").Value = 25114
$ws.Range("N67").Value = -26830

$ws.Range("H101").Value = 12034
$ws.Range("J101").Value = 12034
$ws.Range("L101").Value = 12034
$ws.Range("N101").Value = -18524
